# 24th Oct Sprint 2 Update
# Adds a new Product Backlog item (row 9) on the MAIN sheet and updates the
# status of two existing items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAIN")

# Update existing statuses:
#  - ID 3 (row 7): In Progress -> Complete
#  - ID 4 (row 8): In Progress -> Overdue
$ws.Range("H7").Value = "Complete"
$ws.Range("H8").Value = "Overdue"

# New row 9: "A module to automation tests" backlog item
$ws.Range("C9").Value = "Developer"
$ws.Range("D9").Value = "A module to automation tests"
$ws.Range("E9").Value = "Write testing code and add it to npm commands"
$ws.Range("F9").Value = "Medium"
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = "Overdue"

$ws.Rows("9").RowHeight = 30

# Update selection to match the authored workbook state
$ws.Range("J9").Select()
